# Automatic update of files.
# Bumps Taxonsorteringsordning (column B) from 98931 -> 98932 for every
# matching data row, and swaps the content of rows 11/12 (which describe
# two different species) while also bumping their own B-column values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows whose only change is column B: 98931 -> 98932 ---
$simpleRows = @(2, 3, 4, 5, 6, 7, 8, 9, 10, 13, 14)
foreach ($r in $simpleRows) {
    $ws.Cells.Item($r, 2).Value = 98932
}

# --- Row 11 becomes the "Garnlav / Alectoria sarmentosa" record ---
# (previously on row 12), with its Taxonsorteringsordning bumped from
# 79244 to 79245.
$ws.Range("A11").Value = 130965935
$ws.Range("B11").Value = 79245
$ws.Range("D11").Value = "NT"
$ws.Range("E11").Value = 6425
$ws.Range("F11").Value = "Garnlav"
$ws.Range("G11").Value = "Alectoria sarmentosa"
$ws.Range("H11").Value = "(Ach.) Ach."
$ws.Range("Q11").Value = 496969
$ws.Range("R11").Value = 6713674
$ws.Range("AC11").Value = "Måttlig förekomst . inventering åt vasa vind"
$ws.Range("AX11").Value = "Pia Edfors, Enviro Planning"

# --- Row 12 becomes the "Fläcknycklar / Dactylorhiza maculata" record ---
# (previously on row 11), with its Taxonsorteringsordning bumped from
# 98931 to 98932 (same bump as the rest of the sheet).
$ws.Range("A12").Value = 130965861
$ws.Range("B12").Value = 98932
$ws.Range("D12").Value = "LC"
$ws.Range("E12").Value = 219790
$ws.Range("F12").Value = "Fläcknycklar"
$ws.Range("G12").Value = "Dactylorhiza maculata"
$ws.Range("H12").Value = "(L.) Soó"
$ws.Range("Q12").Value = 497138
$ws.Range("R12").Value = 6713448
$ws.Range("AC12").Value = "Betydelsefulla förekomster . inventering åt vasa vind"
$ws.Range("AX12").Value = "Anders Esplund, Pia Edfors, Enviro Planning"
